# chnages for day 4
# Add a new "InvalidLogin" worksheet after the existing "ValidLogin" sheet,
# populate it with the invalid-login test data, make it the active/selected
# sheet, and clear the "tabSelected" flag from the previously active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet (i.e. after "ValidLogin").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "InvalidLogin"

# Header row.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "password"

# Data row.
$ws.Range("A2").Value = "abcd"
$ws.Range("B2").Value = "xyz"

# Make the new sheet the active tab and set its zoom / selection, matching
# the state captured when the sheet was authored in Excel.
$ws.Activate()
$excel.ActiveWindow.Zoom = 190
$ws.Range("B3").Select()
